# Documentation & Town Update
# - Fill in a few previously-blank achievement cells (C7, C18, D18, F18)
# - Move the saved selection/viewport from C6 to F19 (scrolled down toward row 14)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("This Game is Home Grown" achievement): fill the Name column.
$ws.Range("C7").Value = "This Game is Home Grown"

# Row 18 ("Saskatchewan Shaped, but spooky" achievement): fill Name, Description
# and Picture columns.
$ws.Range("C18").Value = "This Game is Home Grown"
$ws.Range("D18").Value = "Found all collectables"
$ws.Range("F18").Value = "Saskatchewan Shaped, but spooky"

# Update the saved view: scroll so row 14 is at the top and select F19.
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F19").Select()
